# 128. Longest Consecutive Sequence
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match formatting of the row above (row 42) for the new row 43
$ws.Range("A42:D42").Copy()
$ws.Range("A43:D43").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add a new row (row 43) for "Longest Consecutive Sequence"
# (set B before A so shared-string insertion order matches: 56 = Longest
#  Consecutive Sequence, 57 = 128/GFG)
$ws.Range("B43").Value = "Longest Consecutive Sequence"
$ws.Range("A43").Value = "128/GFG"
$ws.Range("C43").Value = "Java"
$ws.Range("D43").Value = 45023

# Update the selected cell as in the saved workbook
$ws.Range("G42").Select()
